# Applies the updated cryptos list values (price + 1h volume change,
# plus a Polygon/WrappedEther row swap) produced by the scheduled scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.962.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +5.19%  "

$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0993"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.090.80"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.837.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.668"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.17%  "

$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.975.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.01%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("E27").Value = "  +3.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.94%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0552"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "

$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.62%  "

$ws.Range("E37").Value = "  +4.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "92.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.339.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0193"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.60%  "

$ws.Range("E41").Value = "  +3.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0670"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.10%  "

